# fix colorful ad info error
$wb = $excel.ActiveWorkbook

# --- Sheet "RTX3090 24G": row 11 is the Colorful iGameRTX3090 Advanced OC entry ---
$ws1 = $wb.Worksheets.Item("RTX3090 24G")
$ws1.Range("D11").Value = "UPI uP9512R?"
$ws1.Range("E11").Value = "UPI uP9512R?"
$ws1.Range("F11").Value = "OnSemi `nNCP302150`n(50A DrMOS)?"

# --- Sheet "RTX3080 10G": row 4 is the Colorful iGameRTX3090 Advanced OC entry ---
$ws2 = $wb.Worksheets.Item("RTX3080 10G")
$ws2.Range("D4").Value = "UPI uP9512R"
$ws2.Range("E4").Value = "UPI uP9512R"
$ws2.Range("F4").Value = "OnSemi `nNCP302150`n(50A DrMOS)"

# --- restore the active selections recorded in the saved workbook ---
$ws1.Activate()
$ws1.Range("F11").Select()

$ws2.Activate()
$ws2.Range("D4:F4").Select()
$ws2.Range("F4").Activate()

# re-activate the first sheet (tabSelected="1" stays on RTX3090 24G)
$ws1.Activate()
